# Update TPM-derived ligand/receptor expression metrics on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Target cluster: ECs)
$ws.Range("G2").Value = 0.135749
$ws.Range("H2").Value = 0.407247
$ws.Range("M2").Value = 8.252454666666667
$ws.Range("N2").Value = 24.757364
$ws.Range("O2").Value = 0.05349680956196952
$ws.Range("P2").Value = 0.05349680956196953
$ws.Range("Q2").Value = 1.120262468545334
$ws.Range("R2").Value = 10.082362216908
$ws.Range("S2").Value = 0.05349680956196952
$ws.Range("T2").Value = 0.05349680956196953

# Row 3 (Target cluster: FAPs)
$ws.Range("G3").Value = 0.135749
$ws.Range("H3").Value = 0.407247
$ws.Range("O3").Value = 0.5638948237978928
$ws.Range("P3").Value = 0.5638948237978929
$ws.Range("Q3").Value = 11.80837161094633
$ws.Range("R3").Value = 106.275344498517
$ws.Range("S3").Value = 0.5638948237978928
$ws.Range("T3").Value = 0.5638948237978929

# Row 4 (Target cluster: MuSCs)
$ws.Range("G4").Value = 0.135749
$ws.Range("H4").Value = 0.407247
$ws.Range("M4").Value = 57.81408433333333
$ws.Range("N4").Value = 173.442253
$ws.Range("O4").Value = 0.3747817085348802
$ws.Range("P4").Value = 0.3747817085348802
$ws.Range("Q4").Value = 7.848204134165667
$ws.Range("R4").Value = 70.633837207491
$ws.Range("S4").Value = 0.3747817085348802
$ws.Range("T4").Value = 0.3747817085348802

# Row 5 (Target cluster: Resolving-Mac)
$ws.Range("G5").Value = 0.135749
$ws.Range("H5").Value = 0.407247
$ws.Range("M5").Value = 1.207345666666667
$ws.Range("N5").Value = 3.622037
$ws.Range("O5").Value = 0.007826658105257385
$ws.Range("P5").Value = 0.007826658105257386
$ws.Range("Q5").Value = 0.1638959669043333
$ws.Range("R5").Value = 1.475063702139
$ws.Range("S5").Value = 0.007826658105257385
$ws.Range("T5").Value = 0.007826658105257386
